$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing the existing row 15 (and below) down to row 16.
$ws.Rows.Item(15).Insert()

# New row 15 data (this is the new record being added)
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C15").Value = "Los Lagos"
$ws.Range("D15").Value = 44491
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 100112013
$ws.Range("G15").Value = "Alcachofa"
$ws.Range("H15").Value = "Madrigal"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = 11000
$ws.Range("N15").Value = "`$/caja 40 unidades"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 275
$ws.Range("Q15").Value = 40
$ws.Range("R15").Value = "Hortaliza"
